$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Science 37" row (row 28), shifting
# that row (and its formulas/hyperlink) down to row 29.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new company (Avantor / AVTR),
# matching the plain "Name/Ticker only" style used by the other rows above it.
$ws.Range("B28").Value = "Avantor"
$ws.Range("C28").Value = "AVTR"

# The row insert shifts the underlying cell data correctly, but the
# worksheet's <hyperlinks> list keeps a stale reference to the old B28
# location. Rebuild the hyperlinks so the SNCE.xlsx link follows its cell
# to B29, while UNH.xlsx/CVS.xlsx stay put on B3/B4.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B29"), "SNCE.xlsx")
$ws.Hyperlinks.Add($ws.Range("B3"), "UNH.xlsx")
$ws.Hyperlinks.Add($ws.Range("B4"), "CVS.xlsx")

# Re-adding the hyperlinks re-applies the "Hyperlink" cell style, which is
# what these cells already had, so this keeps them consistent.
$ws.Range("B29").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"
$ws.Range("B4").Style = "Hyperlink"

# Match the selected cell reflected in the saved file.
$ws.Range("D28").Select()
